$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.61%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.71%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.048"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07945"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.91%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.895"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.04%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9277"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.11%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1347"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "29.94%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1899"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.85%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09095"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.13%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03432"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.59%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09917"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.24%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001392"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.96%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005854"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.50%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.526"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.47%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.031"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.61%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.938"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.26%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.25%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.67%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.053"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.01%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2397"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.35%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04491"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.23%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.48%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004765"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.57%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001231"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.63%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003002"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-32.57%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01890"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.46%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04744"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.81%"
$ws.Range("B41").Value = "Dexo"
$ws.Range("C41").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01043"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "32.35%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007336"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.76%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.40%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01101"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.85%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006284"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.74%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-65.15%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.04%"
